$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell H1 onto I1:J1 so they
# pick up the same bold/border/centered style used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I and J for rows 2-7
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2

$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 8
